$wb = $excel.ActiveWorkbook

# Rename the "RFI" unit label to "a.u." everywhere it is used on the
# Samples sheet (FL1/FL2/FL3 unit columns).
$samples = $wb.Worksheets.Item("Samples")
$samples.Range("E4").Value = "a.u."
$samples.Range("F9").Value = "a.u."
$samples.Range("G10").Value = "a.u."

# Switch the active sheet from "Beads" to "Samples" and move the
# Samples sheet's selection from K11 to G6.
$samples.Activate() | Out-Null
$samples.Range("G6").Select() | Out-Null
